# Split the transmission line from bus 1 to 2 into two separate lines.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Line data")

# Insert a new row right after the existing "1 -> 2" line (row 2),
# shifting all subsequent line records down by one. Copying row 2
# first means the new row inherits its cell formatting (e.g. the
# bus-number style on column A) as Excel would when duplicating it.
$ws.Rows.Item(2).Copy()
$ws.Rows.Item(3).Insert()

# Update the original line (row 2) and the newly inserted duplicate
# (row 3) with the split-line impedance/susceptance values.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 0.038760000000000003
$ws.Range("D2").Value = 0.11834
$ws.Range("E2").Value = 0.0264
$ws.Range("F2").Value = 95

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 0.038760000000000003
$ws.Range("D3").Value = 0.11834
$ws.Range("E3").Value = 0.0264
$ws.Range("F3").Value = 95

# Make "Line data" the active sheet/tab, as it was left selected
# after the edit.
$ws.Activate()
